$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Formula = "=_xlfn.STDEV.S(C2:C3)"
$ws.Range("N3").Formula = "=_xlfn.STDEV.S(D2:D3)"
$ws.Range("N4").Formula = "=_xlfn.STDEV.S(E2:E3)"
$ws.Range("N5").Formula = "=_xlfn.STDEV.S(F2:F3)"

$ws.Range("N6").Select()
